# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across all craft-job sheets with newly fetched
# Universalis pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 555.625
$ws.Range("I4").Value = 492.14285
$ws.Range("K4").Value = 492.14285
$ws.Range("M4").Value = -378.14285

$ws.Range("H15").Value = 1122.4657
$ws.Range("I15").Value = 1122.4657
$ws.Range("K15").Value = 3367.3971
$ws.Range("M15").Value = -3198.3971

$ws.Range("H17").Value = 770888.1
$ws.Range("J17").Value = 770888.1
$ws.Range("L17").Value = 2312664.3
$ws.Range("N17").Value = -2313000.3

$ws.Range("H18").Value = 699.5
$ws.Range("I18").Value = 699.5
$ws.Range("K18").Value = 699.5
$ws.Range("M18").Value = -415.5

$ws.Range("H98").Value = 1698.2778
$ws.Range("I98").Value = 838.8
$ws.Range("J98").Value = 5995.6665
$ws.Range("K98").Value = 838.8
$ws.Range("L98").Value = 5995.6665
$ws.Range("M98").Value = 659.2
$ws.Range("N98").Value = -8991.666499999999

$ws.Range("H112").Value = 1496.2368
$ws.Range("I112").Value = 1107.5
$ws.Range("J112").Value = 1541.9706
$ws.Range("K112").Value = 3322.5
$ws.Range("L112").Value = 4625.9118
$ws.Range("M112").Value = -2214.5
$ws.Range("N112").Value = -6841.9118

$ws.Range("H113").Value = 4683.5713
$ws.Range("I113").Value = 3664.2
$ws.Range("J113").Value = 5249.8887
$ws.Range("K113").Value = 3664.2
$ws.Range("L113").Value = 5249.8887
$ws.Range("M113").Value = -410.1999999999998
$ws.Range("N113").Value = -11757.8887

$ws.Range("H122").Value = 1698.2778
$ws.Range("I122").Value = 838.8
$ws.Range("J122").Value = 5995.6665
$ws.Range("K122").Value = 2516.4
$ws.Range("L122").Value = 17986.9995
$ws.Range("M122").Value = -66.39999999999964
$ws.Range("N122").Value = -22886.9995

$ws.Range("H135").Value = 1543.2273
$ws.Range("I135").Value = 1246.421
$ws.Range("J135").Value = 3423
$ws.Range("K135").Value = 11217.789
$ws.Range("L135").Value = 30807
$ws.Range("M135").Value = -8682.789000000001
$ws.Range("N135").Value = -35877

$ws.Range("H137").Value = 4629.3335
$ws.Range("I137").Value = 1888
$ws.Range("K137").Value = 5664
$ws.Range("M137").Value = -3114

$ws.Range("H141").Value = 1772
$ws.Range("I141").Value = 1191.3334
$ws.Range("K141").Value = 3574.0002
$ws.Range("M141").Value = 1605.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1565.579
$ws.Range("I2").Value = 1403.1875
$ws.Range("K2").Value = 1403.1875
$ws.Range("M2").Value = -1290.1875

$ws.Range("H6").Value = 1000
$ws.Range("J6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("N6").Value = -1346

$ws.Range("H32").Value = 5806.3506
$ws.Range("I32").Value = 4126.422
$ws.Range("K32").Value = 4126.422
$ws.Range("M32").Value = -3839.422

$ws.Range("H45").Value = 9702.6
$ws.Range("I45").Value = 13246.556
$ws.Range("J45").Value = 4386.6665
$ws.Range("K45").Value = 13246.556
$ws.Range("L45").Value = 4386.6665
$ws.Range("M45").Value = -12869.556
$ws.Range("N45").Value = -5140.6665

$ws.Range("H74").Value = 1777.6086
$ws.Range("I74").Value = 953
$ws.Range("J74").Value = 3323.75
$ws.Range("K74").Value = 953
$ws.Range("L74").Value = 3323.75
$ws.Range("M74").Value = -79
$ws.Range("N74").Value = -5071.75

$ws.Range("H77").Value = 1777.6086
$ws.Range("I77").Value = 953
$ws.Range("J77").Value = 3323.75
$ws.Range("K77").Value = 4765
$ws.Range("L77").Value = 16618.75
$ws.Range("M77").Value = -397
$ws.Range("N77").Value = -25354.75

$ws.Range("H97").Value = 791
$ws.Range("J97").Value = 299.66666
$ws.Range("L97").Value = 299.66666
$ws.Range("N97").Value = -1291.66666

$ws.Range("H116").Value = 1565.579
$ws.Range("I116").Value = 1403.1875
$ws.Range("K116").Value = 1403.1875
$ws.Range("M116").Value = 890.8125

$ws.Range("H122").Value = 2607.1738
$ws.Range("I122").Value = 1808.1428
$ws.Range("K122").Value = 5424.428400000001
$ws.Range("M122").Value = -2974.428400000001

$ws.Range("H132").Value = 2044.6346
$ws.Range("I132").Value = 1843.4897
$ws.Range("J132").Value = 5330
$ws.Range("K132").Value = 5530.4691
$ws.Range("L132").Value = 15990
$ws.Range("M132").Value = -3000.4691
$ws.Range("N132").Value = -21050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1565.579
$ws.Range("I3").Value = 1403.1875
$ws.Range("K3").Value = 1403.1875
$ws.Range("M3").Value = -1289.1875

$ws.Range("H20").Value = 26323690
$ws.Range("I20").Value = 31258300
$ws.Range("J20").Value = 5769
$ws.Range("K20").Value = 31258300
$ws.Range("L20").Value = 5769
$ws.Range("M20").Value = -31258053
$ws.Range("N20").Value = -6263

$ws.Range("H86").Value = 2989.8572
$ws.Range("I86").Value = 2425.9
$ws.Range("J86").Value = 4399.75
$ws.Range("K86").Value = 2425.9
$ws.Range("L86").Value = 4399.75
$ws.Range("M86").Value = -1302.9
$ws.Range("N86").Value = -6645.75

$ws.Range("H89").Value = 2989.8572
$ws.Range("I89").Value = 2425.9
$ws.Range("J89").Value = 4399.75
$ws.Range("K89").Value = 12129.5
$ws.Range("L89").Value = 21998.75
$ws.Range("M89").Value = -6513.5
$ws.Range("N89").Value = -33230.75

$ws.Range("H96").Value = 8500
$ws.Range("I96").Value = 8500
$ws.Range("K96").Value = 8500
$ws.Range("M96").Value = -5754

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 935.61536
$ws.Range("I16").Value = 568.75
$ws.Range("K16").Value = 568.75
$ws.Range("M16").Value = -281.75

$ws.Range("H31").Value = 39253.52
$ws.Range("I31").Value = 51506.55
$ws.Range("K31").Value = 51506.55
$ws.Range("M31").Value = -51211.55

$ws.Range("H34").Value = 39253.52
$ws.Range("I34").Value = 51506.55
$ws.Range("K34").Value = 51506.55
$ws.Range("M34").Value = -51304.55

$ws.Range("H86").Value = 3430.6428
$ws.Range("I86").Value = 3198.625
$ws.Range("K86").Value = 3198.625
$ws.Range("M86").Value = -2075.625

$ws.Range("H89").Value = 3430.6428
$ws.Range("I89").Value = 3198.625
$ws.Range("K89").Value = 15993.125
$ws.Range("M89").Value = -10377.125

$ws.Range("H93").Value = 20300.2

$ws.Range("H99").Value = 25838.834
$ws.Range("I99").Value = 45670.332
$ws.Range("K99").Value = 45670.332
$ws.Range("M99").Value = -44172.332

$ws.Range("H107").Value = 1979.125
$ws.Range("I107").Value = 1499.2858
$ws.Range("K107").Value = 1499.2858
$ws.Range("M107").Value = 420.7141999999999

$ws.Range("H113").Value = 935.61536
$ws.Range("I113").Value = 568.75
$ws.Range("K113").Value = 568.75
$ws.Range("M113").Value = 1601.25

$ws.Range("H126").Value = 25838.834
$ws.Range("I126").Value = 45670.332
$ws.Range("K126").Value = 137010.996
$ws.Range("M126").Value = -134540.996

$ws.Range("H132").Value = 4763.3184
$ws.Range("I132").Value = 4977.4443
$ws.Range("K132").Value = 14932.3329
$ws.Range("M132").Value = -12402.3329

$ws.Range("H134").Value = 2454.4443
$ws.Range("J134").Value = 4486.4287
$ws.Range("L134").Value = 13459.2861
$ws.Range("N134").Value = -18529.2861

$ws.Range("H141").Value = 183081.89
$ws.Range("J141").Value = 183081.89
$ws.Range("L141").Value = 183081.89
$ws.Range("N141").Value = -193441.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3032.3333
$ws.Range("I109").Value = 2448.5
$ws.Range("K109").Value = 7345.5
$ws.Range("M109").Value = -6305.5

$ws.Range("H113").Value = 1828.5862
$ws.Range("J113").Value = 1937.5385
$ws.Range("L113").Value = 5812.6155
$ws.Range("N113").Value = -10152.6155

$ws.Range("H122").Value = 782.3570999999999
$ws.Range("J122").Value = 608.1667
$ws.Range("L122").Value = 5473.5003
$ws.Range("N122").Value = -10373.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3099.25
$ws.Range("I126").Value = 3199.5715
$ws.Range("K126").Value = 9598.7145
$ws.Range("M126").Value = -7128.7145

$ws.Range("H132").Value = 4335.6577
$ws.Range("I132").Value = 2909.5
$ws.Range("K132").Value = 8728.5
$ws.Range("M132").Value = -6198.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 46075.7
$ws.Range("I7").Value = 58100.8
$ws.Range("K7").Value = 58100.8
$ws.Range("M7").Value = -57988.8

$ws.Range("H126").Value = 46075.7
$ws.Range("I126").Value = 58100.8
$ws.Range("K126").Value = 174302.4
$ws.Range("M126").Value = -171832.4

$ws.Range("H132").Value = 4488.933
$ws.Range("I132").Value = 3791.76
$ws.Range("K132").Value = 11375.28
$ws.Range("M132").Value = -8845.280000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 9135.5
$ws.Range("J69").Value = 9135.5
$ws.Range("L69").Value = 9135.5
$ws.Range("N69").Value = -10633.5

$ws.Range("H72").Value = 9135.5
$ws.Range("J72").Value = 9135.5
$ws.Range("L72").Value = 27406.5
$ws.Range("N72").Value = -34894.5

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -9530
